$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the country-name row (row 5, columns A-C) ---
# These cells used to repeat the long indicator description; they now hold
# the country name in Kyrgyz, Russian and English.
$ws.Range("A5").Value2 = "Кыргыз Республикасы"
$ws.Range("B5").Value2 = "Кыргызская Республика"
$ws.Range("C5").Value2 = "Kyrgyz Republic "

# --- Merge the widths of columns A:C into a single uniform width ---
$ws.Columns("A:C").ColumnWidth = 35

# --- Add the new year column (2023) ---
# Copy formatting from the previous year's cells so the new cells match style.
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R4").Value2 = 2023
$ws.Range("R5").Value2 = 53.5

# --- Adjust row 5 height ---
$ws.Rows("5:5").RowHeight = 21

# --- Reset the selection to the top-left cell ---
$ws.Range("A1").Select() | Out-Null
